$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 41
$ws.Range("H41").Value = 146.11111
$ws.Range("I41").Value = 175
$ws.Range("J41").Value = 137.85715
$ws.Range("K41").Value = 175
$ws.Range("L41").Value = 137.85715
$ws.Range("M41").Value = 265
$ws.Range("N41").Value = -1017.85715
# Row 112
$ws.Range("H112").Value = 6241.5264
$ws.Range("J112").Value = 6319.0537
$ws.Range("L112").Value = 18957.1611
$ws.Range("N112").Value = -21173.1611
# Row 138
$ws.Range("H138").Value = 2319.575
$ws.Range("J138").Value = 2447.9824
$ws.Range("L138").Value = 7343.9472
$ws.Range("N138").Value = -17623.9472

$ws = $wb.Worksheets.Item("ARM")
# Row 122
$ws.Range("H122").Value = 47187.773
$ws.Range("I122").Value = 60367.766
$ws.Range("K122").Value = 181103.298
$ws.Range("M122").Value = -178653.298
# Row 123
$ws.Range("H123").Value = 26457.4
$ws.Range("J123").Value = 26457.4
$ws.Range("L123").Value = 26457.4
$ws.Range("N123").Value = -36257.4

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 1356.7142
$ws.Range("I94").Value = 832.3333
$ws.Range("J94").Value = 1750
$ws.Range("K94").Value = 832.3333
$ws.Range("L94").Value = 1750
$ws.Range("M94").Value = -381.3333
$ws.Range("N94").Value = -2652
# Row 134
$ws.Range("H134").Value = 3125.72
$ws.Range("I134").Value = 3188
$ws.Range("J134").Value = 3058.25
$ws.Range("K134").Value = 9564
$ws.Range("L134").Value = 9174.75
$ws.Range("M134").Value = -7029
$ws.Range("N134").Value = -14244.75

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 5677.784
$ws.Range("I31").Value = 1375.3684
$ws.Range("J31").Value = 8232.343999999999
$ws.Range("K31").Value = 1375.3684
$ws.Range("L31").Value = 8232.343999999999
$ws.Range("M31").Value = -1080.3684
$ws.Range("N31").Value = -8822.343999999999
# Row 34
$ws.Range("H34").Value = 5677.784
$ws.Range("I34").Value = 1375.3684
$ws.Range("J34").Value = 8232.343999999999
$ws.Range("K34").Value = 1375.3684
$ws.Range("L34").Value = 8232.343999999999
$ws.Range("M34").Value = -1173.3684
$ws.Range("N34").Value = -8636.343999999999
# Row 132
$ws.Range("H132").Value = 9261968
$ws.Range("I132").Value = 3102.4
$ws.Range("J132").Value = 12823070
$ws.Range("K132").Value = 9307.200000000001
$ws.Range("L132").Value = 38469210
$ws.Range("M132").Value = -6777.200000000001
$ws.Range("N132").Value = -38474270
# Row 134
$ws.Range("H134").Value = 1726.7142
$ws.Range("I134").Value = 1554.6666
$ws.Range("K134").Value = 4663.9998
$ws.Range("M134").Value = -2128.9998

$ws = $wb.Worksheets.Item("CUL")
# Row 46
$ws.Range("H46").Value = 639.25
$ws.Range("I46").Value = 639.25
$ws.Range("K46").Value = 1917.75
$ws.Range("M46").Value = -1826.75
# Row 70
$ws.Range("H70").Value = 1221.8334
$ws.Range("I70").Value = 969.2727
$ws.Range("K70").Value = 2907.8181
$ws.Range("M70").Value = -2592.8181
# Row 73
$ws.Range("H73").Value = 1221.8334
$ws.Range("I73").Value = 969.2727
$ws.Range("K73").Value = 2907.8181
$ws.Range("M73").Value = -1815.8181
# Row 100
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
# Row 106
$ws.Range("H106").Value = 7742.5
$ws.Range("J106").Value = 8885.799999999999
$ws.Range("L106").Value = 26657.4
$ws.Range("N106").Value = -28549.4
# Row 112
$ws.Range("H112").Value = 3433.3333
$ws.Range("I112").Value = 3433.3333
$ws.Range("J112").Value = 0
$ws.Range("K112").Value = 10299.9999
$ws.Range("L112").Value = 0
$ws.Range("M112").Value = -9191.999899999999
$ws.Range("N112").ClearContents()
# Row 122
$ws.Range("H122").Value = 6087.3887
$ws.Range("I122").Value = 427.5
$ws.Range("J122").Value = 13162.25
$ws.Range("K122").Value = 3847.5
$ws.Range("L122").Value = 118460.25
$ws.Range("M122").Value = -1397.5
$ws.Range("N122").Value = -123360.25
# Row 123
$ws.Range("H123").Value = 6000
$ws.Range("J123").Value = 8000
$ws.Range("L123").Value = 24000
$ws.Range("N123").Value = -28900
# Row 125
$ws.Range("H125").Value = 1845.3125
$ws.Range("I125").Value = 1000
$ws.Range("J125").Value = 1966.0714
$ws.Range("K125").Value = 3000
$ws.Range("L125").Value = 5898.2142
$ws.Range("M125").Value = 1920
$ws.Range("N125").Value = -15738.2142
# Row 138
$ws.Range("H138").Value = 3916.1428
$ws.Range("J138").Value = 7930
$ws.Range("L138").Value = 23790
$ws.Range("N138").Value = -34070

$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 3006.2334
$ws.Range("I122").Value = 1763.375
$ws.Range("J122").Value = 3458.182
$ws.Range("K122").Value = 5290.125
$ws.Range("L122").Value = 10374.546
$ws.Range("M122").Value = -2840.125
$ws.Range("N122").Value = -15274.546
# Row 126
$ws.Range("H126").Value = 1907.8889
$ws.Range("I126").Value = 1887.7368
$ws.Range("J126").Value = 1955.75
$ws.Range("K126").Value = 5663.2104
$ws.Range("L126").Value = 5867.25
$ws.Range("M126").Value = -3193.2104
$ws.Range("N126").Value = -10807.25
# Row 132
$ws.Range("H132").Value = 3125.1765
$ws.Range("I132").Value = 2851.5
$ws.Range("J132").Value = 3368.4443
$ws.Range("K132").Value = 8554.5
$ws.Range("L132").Value = 10105.3329
$ws.Range("M132").Value = -6024.5
$ws.Range("N132").Value = -15165.3329

$ws = $wb.Worksheets.Item("LTW")
# Row 93
$ws.Range("H93").Value = 11336.4
$ws.Range("I93").Value = 15709.286
$ws.Range("J93").Value = 1133
$ws.Range("K93").Value = 15709.286
$ws.Range("L93").Value = 1133
$ws.Range("M93").Value = -14461.286
$ws.Range("N93").Value = -3629
# Row 100
$ws.Range("H100").Value = 2599.3
$ws.Range("I100").Value = 2399
$ws.Range("J100").Value = 2799.6
$ws.Range("K100").Value = 2399
$ws.Range("L100").Value = 2799.6
$ws.Range("M100").Value = -1858
$ws.Range("N100").Value = -3881.6
# Row 122
$ws.Range("H122").Value = 3126.5881
$ws.Range("I122").Value = 3057.5557
$ws.Range("J122").Value = 3392.8572
$ws.Range("K122").Value = 9172.667099999999
$ws.Range("L122").Value = 10178.5716
$ws.Range("M122").Value = -6722.667099999999
$ws.Range("N122").Value = -15078.5716

$ws = $wb.Worksheets.Item("WVR")
# Row 7
$ws.Range("H7").Value = 52629.25
$ws.Range("I7").Value = 502
$ws.Range("K7").Value = 502
$ws.Range("M7").Value = -389
# Row 96
$ws.Range("H96").Value = 4832.8613
$ws.Range("I96").Value = 3172.875
$ws.Range("J96").Value = 5307.143
$ws.Range("K96").Value = 3172.875
$ws.Range("L96").Value = 5307.143
$ws.Range("M96").Value = -1799.875
$ws.Range("N96").Value = -8053.143
# Row 113
$ws.Range("H113").Value = 1859.1818
$ws.Range("I113").Value = 2205.6667
$ws.Range("J113").Value = 300
$ws.Range("K113").Value = 6617.000100000001
$ws.Range("L113").Value = 900
$ws.Range("M113").Value = -4447.000100000001
$ws.Range("N113").Value = -5240
# Row 122
$ws.Range("H122").Value = 2335.9333
$ws.Range("I122").Value = 1568.4286
$ws.Range("J122").Value = 3007.5
$ws.Range("K122").Value = 4705.2858
$ws.Range("L122").Value = 9022.5
$ws.Range("M122").Value = -2255.2858
$ws.Range("N122").Value = -13922.5
# Row 123
$ws.Range("H123").Value = 21385
$ws.Range("J123").Value = 21385
$ws.Range("L123").Value = 21385
$ws.Range("N123").Value = -31185
# Row 126
$ws.Range("H126").Value = 2283.7693
$ws.Range("I126").Value = 1785.5
$ws.Range("J126").Value = 3081
$ws.Range("K126").Value = 5356.5
$ws.Range("L126").Value = 9243
$ws.Range("M126").Value = -2886.5
$ws.Range("N126").Value = -14183
# Row 136
$ws.Range("H136").Value = 2739.818
$ws.Range("I136").Value = 2628.0476
$ws.Range("K136").Value = 7884.1428
$ws.Range("M136").Value = -5334.1428
Write-Output "Applied all 34 profit-calc row updates across 8 sheets"
